$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "price"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "language"

# Row 2
$ws.Range("A2").Value = "item7"
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = "nachh"
$ws.Range("D2").Value = "np"

# Row 3
$ws.Range("A3").Value = "item4"
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = "dance"
$ws.Range("D3").Value = "en"

# Update selection to match the final state
$ws.Range("D3").Select()
